$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 9565
$ws.Range("B3").Value = "Vitality"
$ws.Range("C3").Value = 55.6
$ws.Range("D3").Value = 59.4
$ws.Range("E3").Value = 62.3
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 68.1
$ws.Range("H3").Value = 62.6
$ws.Range("I3").Value = 67.9
$ws.Range("J3").Value = 71.4
